# Wrap table recognizer: the recognizer now emits an (empty) cell for every
# column of a row it touches, instead of leaving untouched columns absent
# from the sheet. Bring rows 16, 17, 19 and 20 in line with that: make sure
# B/C have a (blank) cell present, and blank out the stray "。" placeholder
# that had leaked into B19.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell that is already an empty/blank cell (used as the copy source so the
# new cells come out as genuinely blank, not just cleared).
$blank = $ws.Range("B2")

# Row 16 ("细胞管型"): add blank B16 / C16.
$blank.Copy($ws.Range("B16"))
$blank.Copy($ws.Range("C16"))

# Row 17 ("颗粒管型"): add blank C17 (B17 keeps its existing "。").
$blank.Copy($ws.Range("C17"))

# Row 19 ("尿酸盐结晶"): B19 held a stray "。" - blank it out, then add C19.
$ws.Range("B19").ClearContents()
$blank.Copy($ws.Range("B19"))
$blank.Copy($ws.Range("C19"))

# Row 20 ("其他结晶"): add blank B20 / C20.
$blank.Copy($ws.Range("B20"))
$blank.Copy($ws.Range("C20"))
